# Applies the diff: inserts two new rows above the existing header row,
# turning the old text header row (row 1) into row 3, shifting all data
# rows down by 2, and populating the two new rows with their new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top; existing rows 1..21 shift down to 3..23.
$ws.Rows.Item(1).Resize(2).Insert()

# New row 1: numeric sequence 0..12 across columns A..M. (Style index 1 -
# the bold/bordered header style - is carried down automatically with the
# shifted cells, so the freshly-inserted row 1 keeps using it too.)
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item(1, $c).Value = ($c - 1)
}

# New row 2: mostly blank, except C2 = "Flange" and F2 = "Drive".
$ws.Range("C2").Value = "Flange"
$ws.Range("F2").Value = "Drive"

# Row 3 (previously row 1) keeps its header text, but the "thread_size" and
# "material_surface" labels that used to live in L3/M3 are removed.
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = ""

$ws.Range("A1").Select()
